# Update the timesheet week-start date on Sheet1.
# The value is stored as an Excel serial date (42772 -> 2/6/2017,
# 42786 -> 2/20/2017). Sheet4!A1 references Sheet1!A1 via formula and
# will recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = Get-Date -Year 2017 -Month 2 -Day 20 -Hour 0 -Minute 0 -Second 0

$excel.CalculateFullRebuild()
